# Apply the audit-SEO.xlsx edits:
#  - A2: "SEO - Accessibilité" -> "SEO - Accessibilité - Performance"
#  - C2: typo fix "affichange" -> "affichage"
#  - B3: typo fix "rensigné" -> "renseigné"
#  - C3: typo fix "améne" -> "amène"
#  - B8: typo fix "Toogle" -> "Toggle"
#  - A9: "SEO" -> "SEO - Performance"
#  - F15: new hyperlink text + link to MDN accessibility checklist
#  - selection ends on C30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SEO - Accessibilité - Performance"
$ws.Range("C2").Value = "site plus volumineux et mauvais contenu SEO, probléme de d'affichage en version mobile"

$ws.Range("B3").Value = "Pas de titre renseigné pour la page d'acceuil et non explicite pour la page 2"
$ws.Range("C3").Value = "aucune information donné dans la balise title amène une mauvaise SEO"

$ws.Range("B8").Value = "Toggle de navigation obsoléte sur la page 2"

$ws.Range("A9").Value = "SEO - Performance"

$ws.Hyperlinks.Add($ws.Range("F15"), "https://developer.mozilla.org/fr/docs/Accessibilit%C3%A9/Checklist_accessibilite_mobile", [Type]::Missing, [Type]::Missing, "https://developer.mozilla.org/fr/docs/Accessibilit%C3%A9/Checklist_accessibilite_mobile")
$ws.Range("F15").Style = "Lien hypertexte"

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

$ws.Range("C30").Select()
